$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text format for cells whose new numeric-looking price values must remain as text
# (otherwise Excel auto-converts strings like "0.9999" or "1.000" into numbers, losing formatting)
$textCells = @("D4","D5","D6","D8","D9","D10","D11","D13","D15","D16","D18","D20","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D41","D42","D44","D45","D46","D49","D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "29.469.71"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.851.99"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "240.91"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "0.6309"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.07687"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "0.2934"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "24.76"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "0.07750"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "1.867.31"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "5.038"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "0.00001072"
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("D16").Value = "83.76"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "2.127.50"
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("D18").Value = "6.202"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "29.486.32"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "229.05"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "7.468"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "157.29"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "0.1386"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "8.415"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "17.71"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "1.336"
$ws.Range("E29").Value = "  +5.59%  "
$ws.Range("D30").Value = "1.467"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "0.05696"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "4.134"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "4.043"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "1.854"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").Value = "1.166"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").Value = "0.7087"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "2.588"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "0.01794"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").Value = "1.220.91"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "6.547"
$ws.Range("E41").Value = "  +5.33%  "
$ws.Range("D42").Value = "0.9082"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "2.035.24"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "101.99"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").Value = "66.61"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "0.4026"
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("D50").Value = "9.047"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  +0.92%  "
